$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark from its original location (right
#    after the "Windows 10 Fall Creators Update SDK (16299)" run). Word
#    regenerates "_GoBack" automatically at the location of the most recent
#    edit, so we will add a fresh one at the end of the document later.
# ---------------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 2. Insert the new "Known issues" section before the "Update history"
#    heading.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Update history") | Out-Null
$rng.Collapse(1) | Out-Null
$rng.InsertParagraphBefore()

# Re-find "Update history" so we can get a handle on the paragraph that was
# just inserted before it (it is now the paragraph immediately preceding).
$rng = $d.Content
$rng.Find.Execute("Update history") | Out-Null
$updateHistoryParaIndex = $rng.Paragraphs(1).Index
$newHeadingPara = $d.Paragraphs($updateHistoryParaIndex - 1)
$newHeadingPara.Range.Text = "Known issues"
$newHeadingPara.Style = "Heading 1"

# Insert: empty paragraph, bullet paragraph, empty paragraph, NOTE paragraph
# all between the new heading and "Update history".
$afterHeading = $newHeadingPara.Range
$afterHeading.Collapse(0) | Out-Null
$afterHeading.InsertParagraphAfter()

$rng = $d.Content
$rng.Find.Execute("Update history") | Out-Null
$updateHistoryParaIndex = $rng.Paragraphs(1).Index
$emptyPara1 = $d.Paragraphs($updateHistoryParaIndex - 1)

$bulletRange = $emptyPara1.Range
$bulletRange.Collapse(0) | Out-Null
$bulletRange.InsertParagraphAfter()

$rng = $d.Content
$rng.Find.Execute("Update history") | Out-Null
$updateHistoryParaIndex = $rng.Paragraphs(1).Index
$bulletPara = $d.Paragraphs($updateHistoryParaIndex - 1)
$bulletPara.Range.Text = [char]0x201C + "placeholder" + [char]0x201D
$bulletPara.Range.Text = "The shaders in this sample don" + [char]0x2019 + "t support an alpha channel"
$bulletPara.Style = "List Paragraph"
$bulletPara.Range.ListFormat.ApplyListTemplate($word.ListGalleries.Item(1).ListTemplates.Item(1))

$afterBullet = $bulletPara.Range
$afterBullet.Collapse(0) | Out-Null
$afterBullet.InsertParagraphAfter()

$rng = $d.Content
$rng.Find.Execute("Update history") | Out-Null
$updateHistoryParaIndex = $rng.Paragraphs(1).Index
$emptyPara2 = $d.Paragraphs($updateHistoryParaIndex - 1)

$noteRange = $emptyPara2.Range
$noteRange.Collapse(0) | Out-Null
$noteRange.InsertParagraphAfter()

$rng = $d.Content
$rng.Find.Execute("Update history") | Out-Null
$updateHistoryParaIndex = $rng.Paragraphs(1).Index
$notePara = $d.Paragraphs($updateHistoryParaIndex - 1)
$notePara.Range.Text = "NOTE: PBREffect and the shaders in this sample have been integrated into the DirectX Tool Kit for DX11PLACEHOLDER / DX12PLACEHOLDER including support for an alpha channel in the albedo texture as well as an optional emissive texture."

# Turn the DX11PLACEHOLDER / DX12PLACEHOLDER text into real hyperlinks.
$rng = $d.Content
$rng.Find.Execute("DX11PLACEHOLDER") | Out-Null
$d.Hyperlinks.Add($rng, "https://github.com/Microsoft/DirectXTK/wiki/PBREffect", "", "", "DX11") | Out-Null

$rng = $d.Content
$rng.Find.Execute("DX12PLACEHOLDER") | Out-Null
$d.Hyperlinks.Add($rng, "https://github.com/Microsoft/DirectXTK12/wiki/PBREffect", "", "", "DX12") | Out-Null

# ---------------------------------------------------------------------------
# 3. Insert a new paragraph about the February 2018 RMA texture channel
#    change, right after "...with HDR output configured." paragraph.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("with HDR output configured.") | Out-Null
$hdrParaIndex = $rng.Paragraphs(1).Index
$hdrPara = $d.Paragraphs($hdrParaIndex)
$r2 = $hdrPara.Range
$r2.Collapse(0) | Out-Null
$r2.InsertParagraphAfter()
$r2.InsertParagraphAfter()

$emptyAfterHdr = $d.Paragraphs($hdrParaIndex + 1)
$rmaPara = $d.Paragraphs($hdrParaIndex + 2)
$rmaPara.Range.Text = "In February 2018, the channel order for the RMA texture was changed to match the glTF2PLACEHOLDER specification: metalnessPLACEHOLDER in blue, roughnessPLACEHOLDER in green, and ambient occlusionPLACEHOLDER in red."

$rng = $d.Content
$rng.Find.Execute("glTF2PLACEHOLDER") | Out-Null
$d.Hyperlinks.Add($rng, "https://github.com/KhronosGroup/glTF/tree/master/specification/2.0", "", "", "glTF2") | Out-Null

$rng = $d.Content
$rng.Find.Execute("metalnessPLACEHOLDER") | Out-Null
$rng.Text = "metalness"
$rng.Italic = 1

$rng = $d.Content
$rng.Find.Execute("roughnessPLACEHOLDER") | Out-Null
$rng.Text = "roughness"
$rng.Italic = 1

$rng = $d.Content
$rng.Find.Execute("ambient occlusionPLACEHOLDER") | Out-Null
$rng.Text = "ambient occlusion"
$rng.Italic = 1

# ---------------------------------------------------------------------------
# 4. Add a fresh "_GoBack" bookmark at the very end of the document, in the
#    final (empty) paragraph right before the section properties.
# ---------------------------------------------------------------------------
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$d.Bookmarks.Add("_GoBack", $endRange) | Out-Null

# ---------------------------------------------------------------------------
# 5. Bump the copyright year shown in the footers from 2017 to 2018.
# ---------------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    foreach ($hf in $sec.Footers) {
        if ($hf.Exists) {
            $hf.Range.Find.Execute("2017", $true, $false, $false, $false, $false, $true, 1, $false, "2018", 2) | Out-Null
        }
    }
}

Write-Output "script complete"
